# Mark the six "Easy" Array problems (rows 7,8,13,14,15 in column C) as
# done (checkmark), and add "not done" (cross) markers for the whole
# "Medium" Array section (rows 19-45 in column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$doneMark = "✔️"
$notDoneMark = "❌"

# Rows in the EASY section that are now marked as done.
$doneRows = @(7, 8, 13, 14, 15)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 3).Value = $doneMark
}

# Rows in the MEDIUM section that now get a "not done" marker in column C.
for ($r = 19; $r -le 45; $r++) {
    $ws.Cells.Item($r, 3).Value = $notDoneMark
}

# Update the view: scroll back to the top, and select C19:C45.
$ws.Range("C19:C45").Select()
